$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 17

$ws.Cells.Item($row, 1).Value = "2025-08-15 13:05:15 UTC"
$ws.Cells.Item($row, 2).Value = "2025-08-15 18:35:15 IST"
$ws.Cells.Item($row, 3).Value = "SKIPPED"
$ws.Cells.Item($row, 4).Value = "No change in PDF. Skipping download & Excel update."
$ws.Cells.Item($row, 5).Value = "https://nalcoindia.com/wp-content/uploads/2025/08/INGOT-15-08-2025.pdf"
$ws.Cells.Item($row, 6).Value = ""
$ws.Cells.Item($row, 7).Value = 0
$ws.Cells.Item($row, 8).Value = ""

$newRowRange = $ws.Range("A17:H17")
$newRowRange.HorizontalAlignment = -4108
$newRowRange.VerticalAlignment = -4108
